# update scripts wuth new tpm
# Recomputed Ligand/Edge expression values (and their derived specificities)
# for the "ECs -> Adm/Calcr" target-cluster group on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs): ligand expression values updated with new TPM data,
# plus the derived specificity / edge-weight values recomputed from them.
$ws.Range("G2").Value = 15.89577633333333
$ws.Range("H2").Value = 47.687329
$ws.Range("I2").Value = 0.286059172443548
$ws.Range("J2").Value = 0.2860591724435479
$ws.Range("Q2").Value = 67.98571081704766
$ws.Range("R2").Value = 611.871397353429
$ws.Range("S2").Value = 0.286059172443548
$ws.Range("T2").Value = 0.2860591724435479

# Row 3 (FAPs): ligand-expression values unchanged, but their
# group-relative specificities shift because row 2 changed.
$ws.Range("I3").Value = 0.6735478078679881
$ws.Range("J3").Value = 0.673547807867988
$ws.Range("S3").Value = 0.6735478078679881
$ws.Range("T3").Value = 0.673547807867988

# Row 4 (MuSCs): same as row 3 - specificities recomputed against the
# new group total; R4 reflects a tiny float re-evaluation as well.
$ws.Range("I4").Value = 0.04039301968846393
$ws.Range("J4").Value = 0.04039301968846393
$ws.Range("R4").Value = 86.39937390919499
$ws.Range("S4").Value = 0.04039301968846393
$ws.Range("T4").Value = 0.04039301968846393
